# Fruta / hortaliza, semanal
# A new weekly price record was added at the top of the data block
# (row 18), pushing all the existing records (rows 18-31) down by one
# row (to rows 19-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18:31 down to 19:32, inserting a blank row 18.
$ws.Rows("18:18").Insert()

# Populate the new row 18 with the latest weekly observation.
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "Terminal La Palmera de La Serena"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44806
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100101
$ws.Range("H18").Value = "Berries"
$ws.Range("I18").Value = 100101001
$ws.Range("J18").Value = "Arándano (blue)"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 15500
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 15750
$ws.Range("Q18").Value = "`$/bandeja 2 kilos"
$ws.Range("R18").Value = "Provincia de Limarí"
$ws.Range("S18").Value = 7875
$ws.Range("T18").Value = 2
